$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row: "Variable"/"Valor" -> "Variable "/"Calificacion"
$ws.Range("A1").Value = "Variable "
$ws.Range("B1").Value = "Calificacion"

# Row labels stay "Prueba 1/2/3" (unchanged text, just re-affirm them)
$ws.Range("A2").Value = "Prueba 1"
$ws.Range("A3").Value = "Prueba 2"
$ws.Range("A4").Value = "Prueba 3"

# Move the active selection to B5, as left by the author after finishing the table
$ws.Range("B5").Select()
